# Commit: "Code to use column names with repalced _ with a space."
#
# The Customers.docx report has a header row whose cells contained the
# raw (Northwind) column identifiers CompanyName / ContactName /
# ContactTitle. The generating code now turns "Column_Name" into
# "Column Name" (underscores -> spaces) before it lands in the template,
# so the header labels read as normal words:
#
#   CompanyName  -> Company Name
#   ContactName  -> Full Name
#   ContactTitle -> Title
#
# Each label is a single, unique run of text in the document, so a
# straightforward whole-word Find/Replace on the document body is
# enough - no need to touch table layout (column widths) or the
# header's background picture, neither of which are semantic content.

$d = $word.ActiveDocument

$d.Content.Find.Execute("CompanyName", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Company Name", 2)

$d.Content.Find.Execute("ContactName", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Full Name", 2)

$d.Content.Find.Execute("ContactTitle", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "Title", 2)
